# Add "Stop" to the dictionary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Append a new row (row 79) with "Stop" in columns A, B and C.
$ws.Range("A79").Value = "Stop"
$ws.Range("B79").Value = "Stop"
$ws.Range("C79").Value = "Stop"

# Update the view state: activate the sheet, scroll so row 73 is at the
# top of the visible area, and select cell B91 (matches the author's
# on-screen state when the edit was made).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.ScrollColumn = 1
$ws.Range("B91").Select()
